$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.102.87'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '1.826.28'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.38'
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6106'
$ws.Range('E6').Value = '  -3.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.009'
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07329'
$ws.Range('E8').Value = '  -2.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2871'
$ws.Range('E9').Value = '  -2.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.72'
$ws.Range('E10').Value = '  -2.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07705'
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('D12').Value = '1.804.66'
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.928'
$ws.Range('E13').Value = '  -1.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6571'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '81.38'
$ws.Range('E15').Value = '  -1.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008889'
$ws.Range('E16').Value = '  -4.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.844'
$ws.Range('E17').Value = '  -2.97%  '
$ws.Range('D18').Value = '29.083.52'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').Value = '2.062.31'
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '235.59'
$ws.Range('E20').Value = '  +5.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.41'
$ws.Range('E21').Value = '  -1.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.010'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.084'
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.34'
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1392'
$ws.Range('E26').Value = '  -0.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.409'
$ws.Range('E27').Value = '  -1.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.60'
$ws.Range('E28').Value = '  -2.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.488'
$ws.Range('E29').Value = '  -1.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05571'
$ws.Range('E30').Value = '  -6.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.065'
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.072'
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.210'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.823'
$ws.Range('E34').Value = '  -1.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7291'
$ws.Range('E35').Value = '  -2.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.131'
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.814'
$ws.Range('E38').Value = '  +1.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01757'
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('D40').Value = '1.197.24'
$ws.Range('E40').Value = '  -2.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.365'
$ws.Range('E41').Value = '  -3.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8829'
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.009'
$ws.Range('E43').Value = '  +0.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.92'
$ws.Range('E44').Value = '  -1.48%  '
$ws.Range('D45').Value = '1.952.71'
$ws.Range('E45').Value = '  -1.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000122'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '64.15'
$ws.Range('E47').Value = '  -2.83%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5127'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3975'
$ws.Range('E49').Value = '  -2.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.973'
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05798'
$ws.Range('E51').Value = '  -0.80%  '
